# "Logged Week 15 and simulated Week 16"
# Appends two more weeks' worth of play-by-play yardage figures to the
# running logs on the YDS and ST sheets, and updates the aggregate
# counting stats on OFF, DEF, ST, TURNS and PEN to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append new per-play yardage numbers to the four running logs
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 2 4 1 7 2 -2 5 3 16 4 6 8 12 6 13 4 -1 6 10 7 3 15 2 8"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 3 43 5 5 4 22 8 4 8 7 9 2 7 1 9 12 5 5 3 5 5 5 3 3 6 12 3"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 4 2 2 4 2 3 7 9 8 -1 4 11 5 4 9 3 5 5 6 0 1 0 5"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 17 6 11 31 7 22 10 3 13 3 24 6 1 16 9 7 6 9 25 14 11 23 -6"

# ---------------------------------------------------------------------
# OFF sheet: updated situational counts for Home (row 2) and Road (row 3)
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 185
$offWs.Range("D2").Value = 11
$offWs.Range("E2").Value = 13
$offWs.Range("F2").Value = 56
$offWs.Range("G2").Value = 64
$offWs.Range("H2").Value = 4
$offWs.Range("J2").Value = 38
$offWs.Range("L2").Value = 253
$offWs.Range("M2").Value = 171
$offWs.Range("O2").Value = 21
$offWs.Range("P2").Value = 14
$offWs.Range("Q2").Value = 512

$offWs.Range("C3").Value = 190
$offWs.Range("D3").Value = 3
$offWs.Range("E3").Value = 27
$offWs.Range("F3").Value = 108
$offWs.Range("G3").Value = 24
$offWs.Range("I3").Value = 57
$offWs.Range("J3").Value = 33
$offWs.Range("N3").Value = 19

# ---------------------------------------------------------------------
# DEF sheet: updated situational counts for Home (row 2) and Road (row 3)
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 150
$defWs.Range("F2").Value = 38
$defWs.Range("G2").Value = 42
$defWs.Range("L2").Value = 239
$defWs.Range("M2").Value = 152
$defWs.Range("Q2").Value = 397

$defWs.Range("B3").Value = 13
$defWs.Range("C3").Value = 153
$defWs.Range("D3").Value = 5
$defWs.Range("E3").Value = 31
$defWs.Range("F3").Value = 82
$defWs.Range("G3").Value = 30
$defWs.Range("H3").Value = 25
$defWs.Range("I3").Value = 58
$defWs.Range("J3").Value = 40
$defWs.Range("N3").Value = 11

# ---------------------------------------------------------------------
# ST sheet: updated kicking/special-teams totals plus appended per-game
# distance logs (RA/RM rows for D and RA/RM columns)
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 73
$stWs.Range("D2").Value = 56
$stWs.Range("F2").Value = 381
$stWs.Range("G2").Value = 377
$stWs.Range("H2").Value = 8
$stWs.Range("J2").Value = 178
$stWs.Range("K2").Value = 175

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 65 63 65 53 56"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 35 21 20 0 9"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 24 26 34 21"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 42"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 0 0"

# ---------------------------------------------------------------------
# TURNS sheet: Road fumbles corrected
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("D3").Value = 5

# ---------------------------------------------------------------------
# PEN sheet: updated penalty counts
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 14
$penWs.Range("D4").Value = 12
$penWs.Range("B5").Value = 2
